{"js": "// Resume content update:\n// 1) \"...server-side implementations, batch processing and client-side rendering.\"\n//    -> \"...server-side implementations/SSR, batch processing and client-side rendering (CSR).\"\n// 2) Joshua C. Martinez's reference phone number\n//    \"+(63) 900 000 0000\" -> \"Direct Line: +(63).54.881.4128\"\n\nconst body = context.document.body;\n\n// --- 1a) insert \"/SSR\" right after \"implementations\" ---\nconst implResults = body.search(\"implementations\", { matchCase: false, matchWholeWord: false });\nimplResults.load(\"items\");\nawait context.sync();\n\nif (implResults.items.length > 0) {\n  implResults.items[0].insertText(\"/SSR\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- 1b) insert \" (CSR)\" right after \"client-side rendering\" ---\nconst csrResults = body.search(\"client-side rendering\", { matchCase: false, matchWholeWord: false });\ncsrResults.load(\"items\");\nawait context.sync();\n\nif (csrResults.items.length > 0) {\n  csrResults.items[0].insertText(\" (CSR)\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- 2) replace the \"+(63) 900 000 0000\" phone number with the direct line ---\nconst phoneResults = body.search(\"+(63) 900 000 0000\", { matchCase: true, matchWholeWord: false });\nphoneResults.load(\"items\");\nawait context.sync();\n\nif (phoneResults.items.length > 0) {\n  phoneResults.items[0].insertText(\"Direct Line: +(63).54.881.4128\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Resume content update:\n# 1) \"...server-side implementations, batch processing and client-side rendering.\"\n#    -> \"...server-side implementations/SSR, batch processing and client-side rendering (CSR).\"\n# 2) Joshua C. Martinez's reference phone number\n#    \"+(63) 900 000 0000\" -> \"Direct Line: +(63).54.881.4128\"\n\n$d = $word.ActiveDocument\n\n# --- 1a) insert \"/SSR\" right after \"implementations\" ---\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Text = \"implementations\"\n$find1.MatchCase = $false\n$find1.MatchWholeWord = $false\n$find1.Forward = $true\n$find1.Wrap = 0\n$found1 = $find1.Execute()\nif ($found1) {\n    $range1.InsertAfter(\"/SSR\")\n}\n\n# --- 1b) insert \" (CSR)\" right after \"client-side rendering\" ---\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"client-side rendering\"\n$find2.MatchCase = $false\n$find2.MatchWholeWord = $false\n$find2.Forward = $true\n$find2.Wrap = 0\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.InsertAfter(\" (CSR)\")\n}\n\n# --- 2) replace the \"+(63) 900 000 0000\" phone number with the direct line ---\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = \"+(63) 900 000 0000\"\n$find3.MatchCase = $true\n$find3.MatchWholeWord = $false\n$find3.Forward = $true\n$find3.Wrap = 0\n$found3 = $find3.Execute()\nif ($found3) {\n    $range3.Text = \"Direct Line: +(63).54.881.4128\"\n}\n"}
